# Updated symbol list on Tue Dec 27 11:41:13 UTC 2022 with GitHub Actions
#
# Applies the latest coinranking.com scrape to the "cryptos" sheet:
#   - refreshed Price (col D) quotes for a batch of already-listed coins
#   - a few coins changed rank, so their whole row (Coin/Link/Price/Volume)
#     shifted to a neighbouring row
#   - a couple of Volume(1h) cells gained/lost their "Bestin24h"/"Worstin24h"
#     trend suffix

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Writing a numeric-looking string via .Value normally gets coerced to a
    # real number (and loses formatting such as trailing zeros). Forcing the
    # cell to Text first keeps it as a string; clearing back to the Normal
    # style afterwards drops the transient "quote prefix" formatting so the
    # cell ends up indistinguishable from a plain text cell.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Straight price refreshes (Price column only)
# ---------------------------------------------------------------------
$priceUpdates = @{
    "D2"  = "243.64"
    "D3"  = "23.05"
    "D6"  = "3.429"
    "D7"  = "6.488"
    "D8"  = "0.8111"
    "D9"  = "0.9268"
    "D10" = "0.1432"
    "D11" = "0.07395"
    "D12" = "0.03181"
    "D14" = "0.09340"
    "D15" = "3.855"
    "D16" = "0.001580"
    "D17" = "0.04708"
    "D40" = "0.03922"
    "D41" = "0.006293"
    "D44" = "0.008361"
    "D45" = "0.00005192"
    "D47" = "0.6704"
}
foreach ($addr in $priceUpdates.Keys) {
    Set-TextValue $addr $priceUpdates[$addr]
}

# ---------------------------------------------------------------------
# Volume(1h) trend-suffix only changes
# ---------------------------------------------------------------------
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

# ---------------------------------------------------------------------
# Rank reshuffle, rows 18-24: "One" jumps up to rank 17 (row 18) and
# everything below it drops down by one row.
# ---------------------------------------------------------------------
$rows18to24 = @(
    @{ Row = 18; B = "One";         C = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one";            D = "0.0005953";  E = "17OneONE" }
    @{ Row = 19; B = "TigerCash";   C = "https://coinranking.com/coin/6hIn06L2+tigercash-tch";           D = "0.005895";   E = "18TigerCashTCH" }
    @{ Row = 20; B = "BitKan";      C = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan";         D = "0.001257";   E = "19BitKanKAN" }
    @{ Row = 21; B = "HotbitToken"; C = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb";     D = "0.004788";   E = "20HotbitTokenHTB" }
    @{ Row = 22; B = "NitroEx";     C = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx";          D = "0.00008004"; E = "21NitroExNTX" }
    @{ Row = 23; B = "LEO";         C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo";             D = "3.557";      E = "22LEOLEO" }
    @{ Row = 24; B = "BTSEToken";   C = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse";      D = "2.133";      E = "23BTSETokenBTSE" }
)
foreach ($entry in $rows18to24) {
    $r = $entry.Row
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    Set-TextValue "D$r" $entry.D
    $ws.Range("E$r").Value = $entry.E
}

# ---------------------------------------------------------------------
# Rank reshuffle, rows 42-43: CEJI and BKEXToken swap places.
# ---------------------------------------------------------------------
$rows42to43 = @(
    @{ Row = 42; B = "CEJI";      C = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji";   D = "0.003502"; E = "41CEJICEJIBestin24h" }
    @{ Row = 43; B = "BKEXToken"; C = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"; D = "0.1073";   E = "42BKEXTokenBKK" }
)
foreach ($entry in $rows42to43) {
    $r = $entry.Row
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    Set-TextValue "D$r" $entry.D
    $ws.Range("E$r").Value = $entry.E
}
